$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "61.025.88"
Set-TextValue "E2" "  +1.21%  "
Set-TextValue "D3" "3.385.93"
Set-TextValue "E3" "  -0.11%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "572.10"
Set-TextValue "E5" "  +0.11%  "
Set-TextValue "D6" "141.56"
Set-TextValue "E6" "  +0.02%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "E8" "  -0.13%  "
Set-TextValue "D9" "7.67"
Set-TextValue "E9" "  +2.61%  "
Set-TextValue "E10" "  -0.85%  "
Set-TextValue "D11" "0.388"
Set-TextValue "E11" "  -1.50%  "
Set-TextValue "D12" "3.966.27"
Set-TextValue "E12" "  -0.04%  "
Set-TextValue "E13" "  +1.97%  "
Set-TextValue "D14" "27.95"
Set-TextValue "E14" "  -1.18%  "
Set-PlainText "B15" "WrappedEther"
Set-PlainText "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D15" "3.399.29"
Set-TextValue "E15" "  +0.27%  "
Set-PlainText "B16" "ShibaInu"
Set-PlainText "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000171"
Set-TextValue "E16" "  +0.12%  "
Set-TextValue "D17" "61.116.18"
Set-TextValue "E17" "  +1.16%  "
Set-TextValue "D18" "6.12"
Set-TextValue "E18" "  -2.76%  "
Set-TextValue "D19" "13.68"
Set-TextValue "E19" "  -3.04%  "
Set-TextValue "D20" "8.96"
Set-TextValue "E20" "  -1.93%  "
Set-TextValue "D21" "384.50"
Set-TextValue "E21" "  -1.20%  "
Set-TextValue "D22" "75.94"
Set-TextValue "E22" "  +3.31%  "
Set-TextValue "E23" "  -1.68%  "
Set-TextValue "D24" "0.999"
Set-TextValue "E24" "  -0.23%  "
Set-TextValue "E25" "  -1.69%  "
Set-TextValue "D26" "0.184"
Set-TextValue "E26" "  +2.72%  "
Set-TextValue "E27" "  -0.06%  "
Set-TextValue "E28" "  -2.76%  "
Set-TextValue "E29" "  -1.07%  "
Set-TextValue "E30" "  -0.42%  "
Set-TextValue "E32" "  -4.05%  "
Set-TextValue "D33" "23.29"
Set-TextValue "E33" "  -2.04%  "
Set-TextValue "D34" "6.97"
Set-TextValue "E34" "  +0.09%  "
Set-TextValue "D35" "166.05"
Set-TextValue "E35" "  -0.95%  "
Set-TextValue "D36" "3.419.84"
Set-TextValue "E36" "  +0.03%  "
Set-TextValue "E37" "  +1.14%  "
Set-TextValue "E38" "  -2.49%  "
Set-TextValue "D40" "26.46"
Set-TextValue "E40" "  -2.40%  "
Set-TextValue "D41" "1.00"
Set-TextValue "E41" "  +0.15%  "
Set-TextValue "E42" "  -0.91%  "
Set-TextValue "E43" "  -2.27%  "
Set-TextValue "E44" "  -2.28%  "
Set-TextValue "E45" "  +0.14%  "
Set-TextValue "D46" "2.459.04"
Set-TextValue "E46" "  -2.87%  "
Set-TextValue "D47" "22.90"
Set-TextValue "E47" "  -1.65%  "
Set-TextValue "D48" "6.66"
Set-TextValue "E48" "  -2.79%  "
Set-PlainText "B49" "dogwifhat"
Set-PlainText "C49" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D49" "2.14"
Set-TextValue "E49" "  +10.03%  "
Set-PlainText "B50" "VeChain"
Set-PlainText "C50" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D50" "0.0263"
Set-TextValue "E50" "  -1.02%  "
Set-TextValue "E51" "  -2.23%  "
